# Update the "想去人数" (F column) figures on both the "展览" and the
# "全部类型" worksheets, which carry duplicated data in this workbook.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F.
$updates = @{
    2  = 160
    3  = 7200
    4  = 5338
    9  = 106
    12 = 198
    13 = 641
    14 = 235
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
